# Hortaliza, Vega Modelo de Temuco - Zanahoria
# Weekly data refresh: insert a new observation row at row 326
# (pushing the existing rows 326-354 down to 327-355) with a new
# "Región del Bíobío" data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 326; this shifts rows
# 326..354 down to 327..355 and keeps the sheet's row 326 style (s="2"
# date format on column D) that Excel copies from the row being pushed.
$ws.Rows.Item(326).Insert()

$ws.Range("A326").Value = 10
$ws.Range("B326").Value = "Vega Modelo de Temuco"
$ws.Range("C326").Value = "La Araucanía"
$ws.Range("D326").Value = 44826
$ws.Range("E326").Value = 9
$ws.Range("F326").Value = 100114013
$ws.Range("G326").Value = "Zanahoria"
$ws.Range("H326").Value = "Sin especificar"
$ws.Range("I326").Value = "Primera"
$ws.Range("J326").Value = 200
$ws.Range("K326").Value = 12000
$ws.Range("L326").Value = 13000
$ws.Range("M326").Value = 12500
$ws.Range("N326").Value = "$/saco 20 kilos"
$ws.Range("O326").Value = "Región del Bíobío"
$ws.Range("P326").Value = 625
$ws.Range("Q326").Value = 20
$ws.Range("R326").Value = "Hortaliza"
